# Update deviation from circle table, error in mean calculation!
#
# The "avg mean" (col E) and "avg std" (col F) values in each of the three
# R-blocks (R=30, R=50, R=70) were recomputed after fixing a bug in the
# mean calculation, and a couple of the raw xoff samples (D19, D20, B20)
# were corrected as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- R = 30 block (rows 2-6) ---
$ws.Range("E2").Value = 1.3
$ws.Range("F2").Value = 0.67

$ws.Range("E3").Value = 1.29
$ws.Range("F3").Value = 0.69

$ws.Range("E4").Value = 1.52
$ws.Range("F4").Value = 0.81

$ws.Range("E5").Value = 1.21
$ws.Range("F5").Value = 0.82

$ws.Range("E6").Value = 1.85
$ws.Range("F6").Value = 1.3

# --- R = 50 block (rows 9-13) ---
$ws.Range("E9").Value = 1.2
$ws.Range("F9").Value = 0.67

$ws.Range("E10").Value = 1.01
$ws.Range("F10").Value = 0.57

$ws.Range("E11").Value = 1.42
$ws.Range("F11").Value = 0.92

$ws.Range("E12").Value = 1.11
$ws.Range("F12").Value = 0.73

$ws.Range("E13").Value = 0.79
$ws.Range("F13").Value = 0.37

# --- R = 70 block (rows 16-20) ---
$ws.Range("E16").Value = 0.9
$ws.Range("F16").Value = 0.45

$ws.Range("E17").Value = 0.81
$ws.Range("F17").Value = 0.43

$ws.Range("E18").Value = 1.06
$ws.Range("F18").Value = 0.65

$ws.Range("D19").Value = -1
$ws.Range("E19").Value = 0.71
$ws.Range("F19").Value = 0.41

$ws.Range("B20").Value = 35
$ws.Range("D20").Value = -7
$ws.Range("E20").Value = 1.69
$ws.Range("F20").Value = 1.23

# --- Update the active cell / selection on the sheet ---
$ws.Range("B20").Select()

# --- Update the window tab ratio (cosmetic UI state) ---
$excel.ActiveWindow.TabRatio = 0.99
